# Applies the KNIVSTA worksheet update:
#  - Reorders data rows 2..115 according to the new scrape order
#    (the underlying source re-sorted/re-fetched records; row identity
#    follows the "Beteckning" (col A) value).
#  - Bumps the "Förändrad" (col C) timestamp from 2026-02-12 (46065) to
#    2026-02-13 (46066) for every data row.
# All other per-row content (values + HYPERLINK formulas in S/T/U/V/W/X/Y/Z)
# travels with its row unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 115
$maxCol = 26   # A..Z
$colC = 3      # "Förändrad" column
$newChanged = 46066

# new_row_index -> old_row_index (both absolute sheet row numbers), in order
# for new rows 2..115
$map = @(2,3,5,4,6,7,8,15,18,10,16,14,11,17,12,13,9,19,20,21,22,23,24,25,26,27,28,29,30,32,34,35,36,38,33,31,37,39,40,57,87,80,61,74,50,99,54,60,49,55,65,64,69,86,66,63,62,70,71,77,72,75,76,73,58,88,89,90,93,42,41,95,53,97,96,98,100,52,59,102,101,56,103,68,78,104,85,107,106,105,108,81,82,110,109,111,112,113,115,114,83,46,84,45,51,44,47,43,79,48,67,91,92,94)

# 1) Snapshot every cell in the data rows before mutating anything.
$snapVals = @{}
$snapIsFormula = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowVals = @{}
    $rowIsF = @{}
    for ($c = 1; $c -le $maxCol; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        if ($cell.HasFormula) {
            $rowVals[$c] = $cell.Formula
            $rowIsF[$c] = $true
        } else {
            $rowVals[$c] = $cell.Value2()
            $rowIsF[$c] = $false
        }
    }
    $snapVals[$r] = $rowVals
    $snapIsFormula[$r] = $rowIsF
}

# 2) Clear the data rows (keeps per-column number formatting: date format on
#    B/C, wrap-text on R) so no stale content lingers from the reorder.
$ws.Range("A" + $firstRow + ":Z" + $lastRow).ClearContents() | Out-Null

# 3) Write each row back at its new position, sourced from the snapshot.
for ($i = 0; $i -lt $map.Length; $i++) {
    $newRow = $firstRow + $i
    $oldRow = $map[$i]
    $rowVals = $snapVals[$oldRow]
    $rowIsF = $snapIsFormula[$oldRow]
    for ($c = 1; $c -le $maxCol; $c++) {
        if ($c -eq $colC) {
            continue
        }
        $val = $rowVals[$c]
        if ($rowIsF[$c]) {
            $ws.Cells.Item($newRow, $c).Formula = $val
        } elseif ($null -ne $val) {
            $ws.Cells.Item($newRow, $c).Value = $val
        }
    }
    # Column C ("Förändrad") always becomes the new changed-date for every row.
    $ws.Cells.Item($newRow, $colC).Value = $newChanged
}

# 4) Re-assert column formatting (ClearContents should have preserved it,
#    but make it explicit/robust).
$ws.Range("B" + $firstRow + ":C" + $lastRow).NumberFormat = "YYYY-MM-DD"
$ws.Range("R" + $firstRow + ":R" + $lastRow).WrapText = $true
